$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Si do ta vlerësonit përzgjedhjen e çështjeve për t’u trajtuar nga organizatat e shoqërisë civile në kuadër të aktiviteteve të tyre?"

$ws.Range("A3").Value = "Të bazuara në analizën e brengave të grupeve që i përfaqësojnë"
$ws.Range("B3").Value = 43

$ws.Range("A4").Value = "Të bazuara në fondet në dispozicion"
$ws.Range("B4").Value = 16

$ws.Range("A5").Value = "Të bazuara në kontaktet e drejtpërdrejta me grupet që i përfaqësojnë"
$ws.Range("B5").Value = 42

$ws.Range("A6:B12").EntireRow.Delete()
